$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dimension/measure tag for "provincia" (F) and "cooperativa" (G)
$ws.Range("F2").Value = "iaest-measure:provincia"
$ws.Range("G2").Value = "iaest-measure:explotaciones-cuyo-titular-es-una-cooperativa-de-produccion"

# Row 3: dim -> medida marker for F and G
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "medida"

# Row 4: datatype for F and G -> xsd:int
$ws.Range("F4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"

# Row 5: remove the mapping file reference for the cooperativa column (G5)
$ws.Range("G5").Clear()

$wb.Save()
